# Auto-generated edit script applying diff to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.649.54"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "2.500.86"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'575.89"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").Value = "'166.86"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("D9").Value = "2.499.30"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +3.16%  "

$ws.Range("D13").Value = "'4.93"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "2.956.77"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "69.556.92"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "2.502.74"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "'11.20"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("E20").Value = "  -4.86%  "

$ws.Range("D21").Value = "'348.05"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'70.39"
$ws.Range("E25").Value = "  +2.65%  "

$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").Value = "'8.76"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("D28").Value = "2.625.60"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").Value = "0.0₃0891"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").Value = "'7.81"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").Value = "'458.99"
$ws.Range("E32").Value = "  -0.64%  "

$ws.Range("E33").Value = "  -3.05%  "

$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").Value = "'157.50"
$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").Value = "'19.07"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "'18.44"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "'0.317"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").Value = "'4.68"
$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("D43").Value = "'1.60"

$ws.Range("D44").Value = "'38.15"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  -4.20%  "

$ws.Range("D46").Value = "'1.09"
$ws.Range("E46").Value = "  -5.84%  "

$ws.Range("D47").Value = "'141.32"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").Value = "'0.519"
$ws.Range("E49").Value = "  -1.84%  "

$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("E51").Value = "  -0.80%  "

